$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Professional summary paragraph - plain text swap.
# "affecting all Black and Asian-American voters" -> "affecting 50M voters"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Data engineering professional with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data engineering professional with 15+ years building systems that matter. Discovered systematic demographic coding errors affecting 50M voters, developed", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: Bullet point under "Partner - Siege Analytics" needs the phrase
# "all Black and Asian-American" replaced by a *bold, colored* "50M" run,
# splitting what is currently one run into three:
#   "... affecting "  |  "50M" (bold, color 2C3E50)  |  " voters, developed ..."
#
# Directly assigning Font properties to a mid-run sub-range does not persist
# in this host, so we build the formatted "50M" run in a disposable scratch
# paragraph at the end of the document, capture it via .FormattedText (which
# DOES splice/format correctly when assigned into another range), paste it
# into place, then remove the scratch paragraph.
# ---------------------------------------------------------------------------
$scratchAnchor = $d.Content
$scratchAnchor.Collapse(0)
$scratchAnchor.InsertParagraphAfter() | Out-Null

$scratchRange = $d.Content
$scratchRange.Collapse(0)
$scratchRange.InsertAfter("50M") | Out-Null
$scratchRange.Font.Bold = $true
$scratchRange.Font.Color = 5258796
$fiftyMFormatted = $scratchRange.FormattedText

$target = $d.Content
$target.Find.Execute("Discovered systematic race coding errors affecting", $true) | Out-Null
$target.Collapse(0)
$target.Find.Execute("all Black and Asian-American", $true) | Out-Null
$destination = $d.Range($target.Start, $target.End)
$destination.FormattedText = $fiftyMFormatted

$d.Paragraphs.Last.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Change 3: Project impact line - plain text swap.
# "affecting all Black and Asian-American voters," -> "affecting 50M voters nationwide,"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved", 2) | Out-Null
